$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3: question + reuse existing answer text from B2
$ws.Range("A3").Value = "Who is animal husbandry Secretary"
$ws.Range("B3").Value = $ws.Range("B2").Value2

# Update selection to match target state
$ws.Range("A16").Select()
